# Apply edits described by the diff for 6.4.1.2.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text fixes (drop trailing period after "6.4.1.2") ---
# B1 = Russian title, C1 = English title (A1 stays the Kyrgyz title, unchanged)
$ws.Range("B1").Value = "6.4.1.2 Потери воды при транспортировке"
$ws.Range("C1").Value = "6.4.1.2 Percentage of water loss during transportation"

# --- Updated data values for year 2022 (column P) ---
$ws.Range("P5").Value = 2388
$ws.Range("P10").Value = 335.3
$ws.Range("P16").Value = 27.3
$ws.Range("P21").Value = 24.3

# --- Update the active selection shown when the sheet is reopened ---
$ws.Range("S3").Select() | Out-Null
